$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artigos Aceitos")

$ws.Range("A7").Value = "Mônica Lima e Souza"
$ws.Range("B7").Value = 2022
$ws.Range("C7").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Cais do Valongo: caminhos da África no Brasil. In: Silvana Terenzi Neuenschwander. (Org.). Patrimônios do Brasil. 1ed.Belo Horizonte: Lucca, 2022, v. , p. 156-165."

$ws.Range("A8").Value = "Mônica Lima e Souza"
$ws.Range("B8").Value = 2022
$ws.Range("C8").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Liberdade, liberdade, abre as asas sobre nós?. Revista Cult, p. 40 - 43, 01 set.  2022."

$ws.Range("A9").Value = "Mônica Lima e Souza"
$ws.Range("B9").Value = 2021
$ws.Range("C9").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. História da África e dos Africanos no Brasil: patrimônios, história pública e reparação. 2021. (Apresentação de Trabalho/Conferência ou palestra)."

$ws.Range("A10").Value = "Mônica Lima e Souza"
$ws.Range("B10").Value = 2021
$ws.Range("C10").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Educação Patrimonial, História Pública e Reparação. 2021. (Apresentação de Trabalho/Conferência ou palestra)."

$ws.Range("A11").Value = "Mônica Lima e Souza"
$ws.Range("B11").Value = 2021
$ws.Range("C11").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica; NASCIMENTO, Patricia ; OLIVEIRA, L. C. . Heranças Africanas no Brasil. 2021. (Apresentação de Trabalho/Conferência ou palestra)."

$ws.Range("A12").Value = "Mônica Lima e Souza"
$ws.Range("B12").Value = 2021
$ws.Range("C12").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Patrimônio afro-brasileiro e políticas públicas: uma questão de direitos. 2021. (Apresentação de Trabalho/Conferência ou palestra)."

$ws.Range("A13").Value = "Mônica Lima e Souza"
$ws.Range("B13").Value = 2021
$ws.Range("C13").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Trajetórias de africanas na luta por liberdade no Brasil. 2021. (Apresentação de Trabalho/Outra)."

$ws.Range("A14").Value = "Mônica Lima e Souza"
$ws.Range("B14").Value = 2021
$ws.Range("C14").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Escravidão e Liberdade: trajetórias africanas no Brasil. 2021. (Apresentação de Trabalho/Outra)."

$ws.Range("A15").Value = "Mônica Lima e Souza"
$ws.Range("B15").Value = 2021
$ws.Range("C15").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Passados sensíveis nos roteiros e aulas de campo na região do Cais do Valongo. 2021. (Apresentação de Trabalho/Outra)."

$ws.Range("A16").Value = "Mônica Lima e Souza"
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Cosmovisão africana. 2021. (Apresentação de Trabalho/Conferência ou palestra)."

$ws.Range("A17").Value = "Mônica Lima e Souza"
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica; LIBERATO, Carlos Franco . Sylviane Diouf: quando os africanos resistem à escravidão. 2022. (Apresentação de Trabalho/Conferência ou palestra)."

$ws.Range("A18").Value = "Mônica Lima e Souza"
$ws.Range("B18").Value = 2022
$ws.Range("C18").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica; NOGUERA, R. . Por uma Educação Antirracista. 2022. (Apresentação de Trabalho/Conferência ou palestra)."

$ws.Range("A19").Value = "Mônica Lima e Souza"
$ws.Range("B19").Value = 2022
$ws.Range("C19").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Educação paras as relações étnicorraciais e o Colégio de Aplicação da UFRJ. 2022. (Apresentação de Trabalho/Conferência ou palestra)."

$ws.Range("A20").Value = "Mônica Lima e Souza"
$ws.Range("B20").Value = 2023
$ws.Range("C20").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Histórias da escravidão e da liberdade na Pequena África carioca. 2023. (Apresentação de Trabalho/Conferência ou palestra)."

$ws.Range("A21").Value = "Mônica Lima e Souza"
$ws.Range("B21").Value = 2023
$ws.Range("C21").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Cais do Valongo: encruzilhada de histórias e memórias negras no Rio de Janeiro. 2023. (Apresentação de Trabalho/Conferência ou palestra)."

$ws.Range("A22").Value = "Mônica Lima e Souza"
$ws.Range("B22").Value = 2023
$ws.Range("C22").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. O Cais do Valongo, a Pequena África e o Arquivo Nacional. 2023. (Apresentação de Trabalho/Outra)."

$ws.Range("A23").Value = "Mônica Lima e Souza"
$ws.Range("B23").Value = 2021
$ws.Range("C23").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica; NOGUERA, R. . Por uma Educação Antirracista.`n`t`t`t`t`t`t2021. (Desenvolvimento de material didático ou instrucional - CadernodeTextoseAtividades)."

$ws.Range("A24").Value = "Mônica Lima e Souza"
$ws.Range("B24").Value = 2021
$ws.Range("C24").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica; BRAGA, Damião ; BRASIL, S. ; MOREIRA, G. M. . MPF assina acordo para valorização do memorial Cais do Valongo, no Rio de Janeiro. 2021.`n`t`t`t`t`t`t`t(Programa de rádio ou TV/Entrevista)."

$ws.Range("A25").Value = "Mônica Lima e Souza"
$ws.Range("B25").Value = 2021
$ws.Range("C25").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica; MUNDURUKU, Daniel . Heranças Africanas no Brasil. 2021."

$ws.Range("A26").Value = "Mônica Lima e Souza"
$ws.Range("B26").Value = 2022
$ws.Range("C26").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica; MONTEIRO, Ana Maria da F Costa ; ROCHA, H. A. B. . Atividade de lançamento: Independência para quem? Espaços de reexistência. 2022.`n`t`t`t`t`t`t`t(Programa de rádio ou TV/Outra)."

$ws.Range("A27").Value = "Mônica Lima e Souza"
$ws.Range("B27").Value = 2022
$ws.Range("C27").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. A história da escravidão resiste aos apagamentos sucessivos na Pequena África carioca. 2022.`n`t`t`t`t`t`t`t(Programa de rádio ou TV/Entrevista)."

$ws.Range("A28").Value = "Mônica Lima e Souza"
$ws.Range("B28").Value = 2022
$ws.Range("C28").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica; SILVA JR, C. . Uma perspectiva negra sobre a independência do Brasil. 2022."

$ws.Range("A29").Value = "Mônica Lima e Souza"
$ws.Range("B29").Value = 2021
$ws.Range("C29").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Nossas Histórias n.48 - Rede de HistoriadorXs NegrXs. 2021."

$ws.Range("A30").Value = "Mônica Lima e Souza"
$ws.Range("B30").Value = 2023
$ws.Range("C30").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. O Museu de História e Cultura Afro-brasileira e a história da população negra no Brasil. 2023.`n`t`t`t`t`t`t`t(Programa de rádio ou TV/Entrevista)."

$ws.Range("A31").Value = "Mônica Lima e Souza"
$ws.Range("B31").Value = 2023
$ws.Range("C31").Value = "SOUZA, Mônica Lima e  OU  LIMA, Mônica OU LIMA, Monica. Vim de lá. 2023."

Write-Output "Rows 7-31 added to Artigos Aceitos sheet"